$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '68.505.42'
$ws.Cells.Item(2, 5).Value = '  +1.04%  '

$ws.Cells.Item(3, 4).Value = '2.652.02'
$ws.Cells.Item(3, 5).Value = '  +1.19%  '

$ws.Cells.Item(4, 5).Value = '  -0.01%  '

$ws.Cells.Item(5, 4).Value = '''599.80'
$ws.Cells.Item(5, 5).Value = '  +0.75%  '

$ws.Cells.Item(6, 4).Value = '''154.93'
$ws.Cells.Item(6, 5).Value = '  +1.42%  '

$ws.Cells.Item(7, 5).Value = '  -0.02%  '

$ws.Cells.Item(8, 4).Value = '''0.548'
$ws.Cells.Item(8, 5).Value = '  +0.70%  '

$ws.Cells.Item(9, 4).Value = '2.650.94'
$ws.Cells.Item(9, 5).Value = '  +1.20%  '

$ws.Cells.Item(10, 5).Value = '  +9.70%  '

$ws.Cells.Item(12, 5).Value = '  +1.41%  '

$ws.Cells.Item(13, 4).Value = '''0.356'
$ws.Cells.Item(13, 5).Value = '  +2.40%  '

$ws.Cells.Item(14, 4).Value = '''28.34'
$ws.Cells.Item(14, 5).Value = '  +2.87%  '

$ws.Cells.Item(15, 5).Value = '  +3.13%  '

$ws.Cells.Item(16, 4).Value = '3.130.57'

$ws.Cells.Item(17, 4).Value = '68.392.28'
$ws.Cells.Item(17, 5).Value = '  +1.02%  '

$ws.Cells.Item(18, 4).Value = '2.659.37'
$ws.Cells.Item(18, 5).Value = '  +1.64%  '

$ws.Cells.Item(19, 4).Value = '''11.53'

$ws.Cells.Item(20, 4).Value = '''367.15'
$ws.Cells.Item(20, 5).Value = '  -1.40%  '

$ws.Cells.Item(21, 4).Value = '''7.53'
$ws.Cells.Item(21, 5).Value = '  +1.48%  '

$ws.Cells.Item(22, 4).Value = '''4.43'
$ws.Cells.Item(22, 5).Value = '  +4.81%  '

$ws.Cells.Item(23, 4).Value = '''4.92'
$ws.Cells.Item(23, 5).Value = '  +2.46%  '

$ws.Cells.Item(24, 5).Value = '  +2.03%  '

$ws.Cells.Item(25, 4).Value = '''73.85'
$ws.Cells.Item(25, 5).Value = '  +1.86%  '

$ws.Cells.Item(26, 5).Value = '  +0.07%  '

$ws.Cells.Item(27, 4).Value = '''9.93'
$ws.Cells.Item(27, 5).Value = '  +0.60%  '

$ws.Cells.Item(28, 4).Value = '''0.0000108'
$ws.Cells.Item(28, 5).Value = '  +4.69%  '

$ws.Cells.Item(29, 4).Value = '2.780.05'
$ws.Cells.Item(29, 5).Value = '  +0.82%  '

# Row 30: now Binance-PegBSC-USD (was Bittensor)
$ws.Cells.Item(30, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(30, 4).Value = '''1.00'
$ws.Cells.Item(30, 5).Value = '  +0.19%  '

# Row 31: now Bittensor (was Binance-PegBSC-USD)
$ws.Cells.Item(31, 2).Value = 'Bittensor'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(31, 4).Value = '''580.29'
$ws.Cells.Item(31, 5).Value = '  -2.07%  '

$ws.Cells.Item(32, 4).Value = '''8.24'
$ws.Cells.Item(32, 5).Value = '  +5.64%  '

$ws.Cells.Item(33, 5).Value = '  +4.86%  '

$ws.Cells.Item(34, 5).Value = '  +1.93%  '

$ws.Cells.Item(35, 4).Value = '''0.132'
$ws.Cells.Item(35, 5).Value = '  +5.09%  '

$ws.Cells.Item(36, 4).Value = '''1.61'
$ws.Cells.Item(36, 5).Value = '  +5.92%  '

$ws.Cells.Item(37, 5).Value = '  -0.01%  '

$ws.Cells.Item(38, 4).Value = '''159.60'
$ws.Cells.Item(38, 5).Value = '  +0.79%  '

$ws.Cells.Item(39, 4).Value = '''19.53'
$ws.Cells.Item(39, 5).Value = '  +2.16%  '

$ws.Cells.Item(40, 5).Value = '  +2.39%  '

$ws.Cells.Item(41, 5).Value = '  +0.77%  '

$ws.Cells.Item(42, 4).Value = '''5.44'
$ws.Cells.Item(42, 5).Value = '  +3.41%  '

$ws.Cells.Item(43, 4).Value = '''2.71'
$ws.Cells.Item(43, 5).Value = '  +0.96%  '

$ws.Cells.Item(44, 4).Value = '0.0₆0334'
$ws.Cells.Item(44, 5).Value = '  +12.59%  '

$ws.Cells.Item(45, 4).Value = '''17.72'
$ws.Cells.Item(45, 5).Value = '  +3.54%  '

$ws.Cells.Item(47, 4).Value = '''40.52'
$ws.Cells.Item(47, 5).Value = '  +0.21%  '

$ws.Cells.Item(48, 4).Value = '''158.08'
$ws.Cells.Item(48, 5).Value = '  +1.16%  '

$ws.Cells.Item(49, 5).Value = '  +3.50%  '

$ws.Cells.Item(50, 5).Value = '  +2.41%  '

$ws.Cells.Item(51, 4).Value = '''22.06'
$ws.Cells.Item(51, 5).Value = '  +3.70%  '
